$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: per-language status text changes from
#     "Ready for handoff" to "Handback transform failed"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column also reflects the new status, and the
#     handback/handoff filename mismatch is reported in the Error Detail
#     column (P) for the 31e50a06 entry (row 3); widen that column.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = "Handback file name: lyrdy0kg.uah is different with handoff file name: 31e50a06-1fd4-48cc-a03e-04b3ecfcde09.9a41a29d22a9aa44f613bd28075c13a5582084d4.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: same status/error reporting for the de-de target.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = "Handback file name: lyrdy0kg.uah is different with handoff file name: 31e50a06-1fd4-48cc-a03e-04b3ecfcde09.9a41a29d22a9aa44f613bd28075c13a5582084d4.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
